# Apply the geometry/roof-slope update to the 3-storey model geometry sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteFormats = -4122

# --- Core input changes: overall building length/width ---
$ws.Range("C3").Value = 16.5
$ws.Range("C4").Value = 12

# --- Give B9/C9 the same bold styling the rest of the summary rows use
#     (copy the format straight from D9, which already has it). ---
$ws.Range("D9").Copy()
$ws.Range("B9").PasteSpecial($xlPasteFormats)
$ws.Range("C9").PasteSpecial($xlPasteFormats)

# --- Re-layout the roof-slope block: shift the existing 8/12 slope row
#     from row 14/16 down to row 14/15, and add a brand-new 3/12 slope
#     entry at rows 17/18. ---

# Grab the formats we need to reuse before we start clearing things out.
$ws.Range("G14").Copy()
$ws.Range("Z90").PasteSpecial($xlPasteFormats)      # quotePrefix text style (->25)
$ws.Range("H14").Copy()
$ws.Range("Z91").PasteSpecial($xlPasteFormats)      # bold, 0.0000-less numeric style (->24)
$ws.Range("H16").Copy()
$ws.Range("Z92").PasteSpecial($xlPasteFormats)      # bold label style (->23)
$ws.Range("I16").Copy()
$ws.Range("Z93").PasteSpecial($xlPasteFormats)      # bold style (->13)

# Clear the old F16:I16 block entirely (content + formatting) -- its
# content moves (with an updated formula) to F15:I15.
$ws.Range("F16:I16").Clear()
$ws.Range("G15:H15").Clear()

# New F15:I15 -- roof height for the existing 8/12 slope, moved to row 15.
$ws.Range("F15").Value = ' <Opt-Roof-Height>'
$ws.Range("Z90").Copy()
$ws.Range("G15").PasteSpecial($xlPasteFormats)
$ws.Range("G15").Value = "'8/12"
$ws.Range("Z92").Copy()
$ws.Range("H15").PasteSpecial($xlPasteFormats)
$ws.Range("H15").Formula = '=10.363+($C$3*H14)'
$ws.Range("Z93").Copy()
$ws.Range("I15").PasteSpecial($xlPasteFormats)
$ws.Range("I15").Value = "m"

# New F17:H17 -- second roof slope (3/12).
$ws.Range("F17").Value = "roof-slope"
$ws.Range("Z90").Copy()
$ws.Range("G17").PasteSpecial($xlPasteFormats)
$ws.Range("G17").Value = "'3/12"
$ws.Range("Z91").Copy()
$ws.Range("H17").PasteSpecial($xlPasteFormats)
$ws.Range("H17").Formula = '=3/12'

# New F18:I18 -- roof height for the new 3/12 slope.
$ws.Range("F18").Value = ' <Opt-Roof-Height>'
$ws.Range("Z90").Copy()
$ws.Range("G18").PasteSpecial($xlPasteFormats)
$ws.Range("G18").Value = "'3/12"
$ws.Range("Z92").Copy()
$ws.Range("H18").PasteSpecial($xlPasteFormats)
$ws.Range("H18").Formula = '=10.363+($C$3*H17)'
$ws.Range("Z93").Copy()
$ws.Range("I18").PasteSpecial($xlPasteFormats)
$ws.Range("I18").Value = "m"

# Clean up the scratch cells used to stash copied formats.
$ws.Range("Z90:Z93").Clear()
$excel.CutCopyMode = 0

# Column F is now wider to fit the " <Opt-Roof-Height>" label.
$ws.Columns("F:F").ColumnWidth = 17.59

# Selection left where the author's cursor ended up.
$ws.Range("F18").Select()
